$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the AR batch file name references (rows 2-5)
$ws.Range("A2:A5").Value = "install_zone_fileQA4AR_Oci2.bat"

# Update the NL batch file name references (rows 6-9)
$ws.Range("A6:A9").Value = "install_zone_fileQA4NL_Oci2.bat"

# Move the active selection to A3
$ws.Activate()
$ws.Range("A3").Select()
